# Update the "Förändrad" (Changed) date column (C) for rows 2-5
# from 2023-09-14 (45183) to 2023-09-15 (45184).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2").Value = 45184
$ws.Range("C3").Value = 45184
$ws.Range("C4").Value = 45184
$ws.Range("C5").Value = 45184
